$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the current row 164,
# pushing the existing rows 164-174 down to 165-175.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record's data.
$ws.Range("A164").Value = 8
$ws.Range("B164").Value = "Terminal La Palmera de La Serena"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 44578
$ws.Range("E164").Value = 4
$ws.Range("F164").Value = 100112021
$ws.Range("G164").Value = "Ají"
$ws.Range("H164").Value = "Americana (o)"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 600
$ws.Range("K164").Value = 13000
$ws.Range("L164").Value = 14000
$ws.Range("M164").Value = 13500
$ws.Range("N164").Value = "$/caja 15 kilos"
$ws.Range("O164").Value = "Provincia de Limarí"
$ws.Range("P164").Value = 900
$ws.Range("Q164").Value = 15
$ws.Range("R164").Value = "Hortaliza"
